# Add data to excel
# - Rename the "外送咖啡服務" amenity column (L) header to "冰塊杯販售"
# - Append 5 new store rows (rows 2-6) with name/address and amenity flags

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename column L header text.
$ws.Range("L1").Value = "冰塊杯販售"

# Row 2: 八斗子
$ws.Range("C2").Value = "八斗子"
$ws.Range("D2").Value = "基隆市中正區北寧路327號1樓"
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 1
$ws.Range("L2").Value = 1
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 1
$ws.Range("S2").Value = 1
$ws.Range("AB2").Value = 1
$ws.Range("AC2").Value = 1
$ws.Range("AD2").Value = 1
$ws.Range("AF2").Value = 1
$ws.Range("AI2").Value = 1

# Row 3: 北寧
$ws.Range("C3").Value = "北寧"
$ws.Range("D3").Value = "基隆市中正區北寧路382號382-5號"
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 1
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = 1
$ws.Range("S3").Value = 1
$ws.Range("AB3").Value = 1
$ws.Range("AC3").Value = 1
$ws.Range("AD3").Value = 1
$ws.Range("AI3").Value = 1

# Row 4: 正濱
$ws.Range("C4").Value = "正濱"
$ws.Range("D4").Value = "基隆市中正區豐稔街27號29號"
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 1
$ws.Range("L4").Value = 1
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = 1
$ws.Range("T4").Value = 1
$ws.Range("AB4").Value = 1
$ws.Range("AC4").Value = 1
$ws.Range("AF4").Value = 1
$ws.Range("AI4").Value = 1

# Row 5: 旭東
$ws.Range("C5").Value = "旭東"
$ws.Range("D5").Value = "基隆市中正區中正路54號"
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 1
$ws.Range("L5").Value = 1
$ws.Range("Q5").Value = 1
$ws.Range("S5").Value = 1
$ws.Range("AB5").Value = 1
$ws.Range("AC5").Value = 1
$ws.Range("AD5").Value = 1
$ws.Range("AI5").Value = 1

# Row 6: 和平島
$ws.Range("C6").Value = "和平島"
$ws.Range("D6").Value = "基隆市中正區和一路125號127號"
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 1
$ws.Range("L6").Value = 1
$ws.Range("Q6").Value = 1
$ws.Range("S6").Value = 1
$ws.Range("T6").Value = 1
$ws.Range("AB6").Value = 1
$ws.Range("AC6").Value = 1
$ws.Range("AD6").Value = 1
$ws.Range("AI6").Value = 1
